# CiAn: updated list of materials
#
# Updates the "NecesarMateriale" (materials needed) sheet: refreshes several
# item descriptions, bumps some quantities to 10, marks two superseded
# TargetLink-related line items with strikethrough, and appends a new
# "Osciloscoape" line. Also refreshes the active sheet/selection so
# "NecesarMateriale" becomes the active tab (matching the author's last
# on-screen position) instead of "OpenTopics".

$wb = $excel.ActiveWorkbook

$wsMateriale = $wb.Worksheets.Item("NecesarMateriale")
$wsCourses   = $wb.Worksheets.Item("Courses Structure")
$wsOpen      = $wb.Worksheets.Item("OpenTopics")

# --- NecesarMateriale: text + quantity updates -----------------------------

$wsMateriale.Range("B2").Value = "PC-uri cu Windows 7 / Windows 10"

$wsMateriale.Range("B4").Value = "Optional: Simscape with SimElectronics and SimPower Systems, Signal Processing"

# Licente TargetLink / Placute TL pentru code deployment are now struck
# through (superseded), and their quantities bumped to 10.
$wsMateriale.Range("B5").Font.Strikethrough = $true
$wsMateriale.Range("C5").Value = 10

$wsMateriale.Range("B6").Font.Strikethrough = $true
$wsMateriale.Range("C6").Value = 10

$wsMateriale.Range("B7").Value = "Placute pentru Embedded Coder (NUCLEO-L496ZG)"
$wsMateriale.Range("C7").Value = 10

$wsMateriale.Range("B8").Value = "Motor (to be decided type and invertor if needed)"
$wsMateriale.Range("C8").Value = 10

$wsMateriale.Range("B9").Value = "Sursa tensiune (minim 20A, 20V)"
$wsMateriale.Range("C9").Value = 10

# New row for the oscilloscopes.
$wsMateriale.Range("B10").Value = "Osciloscoape"
$wsMateriale.Range("C10").Value = 10

# --- View/selection updates -------------------------------------------------

# "Courses Structure": keep same selected cell column/row semantics, just
# move the cursor to where the author left it.
$wsCourses.Range("J22").Select()

# "OpenTopics": move the cursor; this sheet stops being the active tab.
$wsOpen.Range("E18").Select()

# "NecesarMateriale" becomes the active sheet/tab with its own selection.
$wsMateriale.Activate()
$wsMateriale.Range("B27").Select()
